$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.965.19"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.298.00"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.90"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.66"
$ws.Range("E6").Value = "  +4.85%  "
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.09"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.09"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.83"
$ws.Range("E14").Value = "  +15.12%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.655.02"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.321.64"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("E18").Value = "  +4.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.883.31"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.33"
$ws.Range("E20").Value = "  +8.25%  "
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.71"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.22"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("E25").Value = "  +13.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.74"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.28"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.62"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.13"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.56"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.87"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.103"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("E40").Value = "  +4.29%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.36"
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.992.02"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0284"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.01"
$ws.Range("E46").Value = "  +4.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.52"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.84"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.77"
$ws.Range("E49").Value = "  +9.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.526.21"
$ws.Range("E51").Value = "  +3.13%  "
